# Apply crypto price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value2 = '62.978.88'
$ws.Cells.Item(2, 5).Value2 = '  -0.45%  '
# Row 3
$ws.Cells.Item(3, 4).Value2 = '2.471.00'
$ws.Cells.Item(3, 5).Value2 = '  -0.47%  '
# Row 4
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = '@'
$c.Value2 = '1.00'
$c.Style = 'Normal'
$ws.Cells.Item(4, 5).Value2 = '  +0.05%  '
# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value2 = '572.16'
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value2 = '  -0.94%  '
# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value2 = '148.44'
$c.Style = 'Normal'
$ws.Cells.Item(6, 5).Value2 = '  +1.05%  '
# Row 7
$ws.Cells.Item(7, 5).Value2 = '  -0.06%  '
# Row 8
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = '@'
$c.Value2 = '0.530'
$c.Style = 'Normal'
$ws.Cells.Item(8, 5).Value2 = '  -1.76%  '
# Row 9
$ws.Cells.Item(9, 5).Value2 = '  -0.16%  '
# Row 10
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value2 = '0.162'
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Value2 = '  -0.42%  '
# Row 11
$ws.Cells.Item(11, 5).Value2 = '  -1.06%  '
# Row 12
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value2 = '0.349'
$c.Style = 'Normal'
$ws.Cells.Item(12, 5).Value2 = '  -1.27%  '
# Row 13
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value2 = '28.97'
$c.Style = 'Normal'
$ws.Cells.Item(13, 5).Value2 = '  +1.16%  '
# Row 14
$ws.Cells.Item(14, 5).Value2 = '  -1.97%  '
# Row 15
$ws.Cells.Item(15, 4).Value2 = '2.921.18'
$ws.Cells.Item(15, 5).Value2 = '  -0.37%  '
# Row 16
$ws.Cells.Item(16, 4).Value2 = '62.903.07'
$ws.Cells.Item(16, 5).Value2 = '  -0.30%  '
# Row 17
$ws.Cells.Item(17, 4).Value2 = '2.480.99'
$ws.Cells.Item(17, 5).Value2 = '  -0.04%  '
# Row 18
$ws.Cells.Item(18, 5).Value2 = '  -6.69%  '
# Row 19
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value2 = '10.79'
$c.Style = 'Normal'
$ws.Cells.Item(19, 5).Value2 = '  -2.35%  '
# Row 20
$ws.Cells.Item(20, 5).Value2 = '  +3.61%  '
# Row 21
$ws.Cells.Item(21, 5).Value2 = '  +0.59%  '
# Row 22
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value2 = '322.59'
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Value2 = '  -2.20%  '
# Row 23
$ws.Cells.Item(23, 5).Value2 = '  +0.03%  '
# Row 24
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value2 = '10.17'
$c.Style = 'Normal'
$ws.Cells.Item(24, 5).Value2 = '  +4.02%  '
# Row 25
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value2 = '64.95'
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value2 = '  -2.05%  '
# Row 26
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value2 = '653.86'
$c.Style = 'Normal'
$ws.Cells.Item(26, 5).Value2 = '  -2.69%  '
# Row 27
$ws.Cells.Item(27, 5).Value2 = '  -1.29%  '
# Row 28
$ws.Cells.Item(28, 4).Value2 = '0.0₃0973'
$ws.Cells.Item(28, 5).Value2 = '  -2.70%  '
# Row 29
$ws.Cells.Item(29, 5).Value2 = '  +0.63%  '
# Row 30
$ws.Cells.Item(30, 5).Value2 = '  -3.07%  '
# Row 31
$ws.Cells.Item(31, 5).Value2 = '  -2.26%  '
# Row 32
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Value2 = '1.83'
$c.Style = 'Normal'
$ws.Cells.Item(32, 5).Value2 = '  -2.23%  '
# Row 33
$ws.Cells.Item(33, 5).Value2 = '  +0.33%  '
# Row 34
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Value2 = '0.998'
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Value2 = '  -0.02%  '
# Row 35
$ws.Cells.Item(35, 5).Value2 = '  -3.42%  '
# Row 36
$ws.Cells.Item(36, 5).Value2 = '  -2.13%  '
# Row 37
$ws.Cells.Item(37, 5).Value2 = '  -1.37%  '
# Row 38
$ws.Cells.Item(38, 5).Value2 = '  -1.70%  '
# Row 39
$ws.Cells.Item(39, 2).Value2 = 'Monero'
$ws.Cells.Item(39, 3).Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value2 = '150.11'
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Value2 = '  -0.49%  '
# Row 40
$ws.Cells.Item(40, 2).Value2 = 'EthereumClassic'
$ws.Cells.Item(40, 3).Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value2 = '18.55'
$c.Style = 'Normal'
$ws.Cells.Item(40, 5).Value2 = '  -1.35%  '
# Row 41
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value2 = '2.70'
$c.Style = 'Normal'
$ws.Cells.Item(41, 5).Value2 = '  -0.40%  '
# Row 42
$ws.Cells.Item(42, 5).Value2 = '  -1.86%  '
# Row 43
$ws.Cells.Item(43, 4).Value2 = '0.0₆0309'
$ws.Cells.Item(43, 5).Value2 = '  -1.37%  '
# Row 44
$ws.Cells.Item(44, 5).Value2 = '  +0.00%  '
# Row 45
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value2 = '153.31'
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value2 = '  -1.73%  '
# Row 46
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value2 = '15.41'
$c.Style = 'Normal'
$ws.Cells.Item(46, 5).Value2 = '  +1.83%  '
# Row 47
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value2 = '3.56'
$c.Style = 'Normal'
$ws.Cells.Item(47, 5).Value2 = '  -1.37%  '
# Row 48
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value2 = '20.38'
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value2 = '  -0.88%  '
# Row 49
$ws.Cells.Item(49, 5).Value2 = '  -0.17%  '
# Row 50
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value2 = '0.0512'
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value2 = '  -0.54%  '
# Row 51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value2 = '0.0905'
$c.Style = 'Normal'
$ws.Cells.Item(51, 5).Value2 = '  -1.49%  '
